$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Characteristics")

# Insert 10 new rows before the current row 17 ("N patients"), pushing the
# existing "N patients" / "Notes" rows down to 27/28.
$ws.Rows("17:26").Insert()

# New median characteristic rows
$ws.Range("A17").Value = "Median Characteristic 1"
$ws.Range("B17").Value = 0.24186

$ws.Range("A18").Value = "Median Characteristic 2"
$ws.Range("B18").Value = 0.2514

$ws.Range("A19").Value = "Median Characteristic 3"
$ws.Range("B19").Value = 0.25455

$ws.Range("A20").Value = "Median Characteristic 4"
$ws.Range("B20").Value = 0.24272
$ws.Range("B20:L20").HorizontalAlignment = -4108
$ws.Range("B20:L20").VerticalAlignment = -4160
$ws.Range("K20").WrapText = $true

$ws.Range("A21").Value = "Median Characteristic 5"
$ws.Range("B21").Value = 0.2517

$ws.Range("A22").Value = "Median Characteristic 6"
$ws.Range("B22").Value = 0.25333

# New N characteristic rows (proportions)
$ws.Range("A23").Value = "N Characteristic 1 - No"
$ws.Range("B23").Value = 694

$ws.Range("A24").Value = "N Characteristic 1 - Yes"
$ws.Range("B24").Value = 306

$ws.Range("A25").Value = "N Characteristic 2 - No"
$ws.Range("B25").Value = 670

$ws.Range("A26").Value = "N Characteristic 2 - Yes"
$ws.Range("B26").Value = 330

# Keep the active selection where the author left it
$ws.Range("B27").Select() | Out-Null
